# Verify_40V_On_Addition_Deletion_Of_SlotCards.xlsx
# "Updated test data as per new implementation"
#
# The "Loading Details Name" column (F8:F9) used the shared label "40V (A)".
# The new implementation renames that loading-details label to "40V Rail(A)".
# Both data rows reference the same label, so updating them together keeps
# the workbook's shared-string table collapsed to a single entry (no stray
# duplicate string is introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8:F9").Value = "40V Rail(A)"

# Leave the selection where the edit finished (matches the saved cursor
# position after updating the last data row).
$ws.Range("F9").Select()
